$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for I1 and J1 (copy format from H1 which already has the header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:J77
$data = @{
    2 = @(5, 6)
    3 = @(8, 8)
    4 = @(9, 9)
    5 = @(8, 8)
    6 = @(9, 9)
    7 = @(8, 9)
    8 = @(9, 9)
    9 = @(8, 8)
    10 = @(7, 7)
    11 = @(9, 9)
    12 = @(7, 7)
    13 = @(7, 8)
    14 = @(8, 8)
    15 = @(7, 7)
    16 = @(7, 8)
    17 = @(8, 8)
    18 = @(7, 7)
    19 = @(7, 7)
    20 = @(8, 8)
    21 = @(7, 7)
    22 = @(8, 8)
    23 = @(7, 8)
    24 = @(7, 8)
    25 = @(8, 8)
    26 = @(7, 8)
    27 = @(6, 7)
    28 = @(7, 7)
    29 = @(7, 7)
    30 = @(7, 8)
    31 = @(9, 9)
    32 = @(7, 7)
    33 = @(6, 7)
    34 = @(7, 7)
    35 = @(7, 7)
    36 = @(6, 6)
    37 = @(8, 8)
    38 = @(7, 7)
    39 = @(8, 8)
    40 = @(7, 7)
    41 = @(6, 6)
    42 = @(8, 8)
    43 = @(7, 7)
    44 = @(9, 9)
    45 = @(7, 7)
    46 = @(5, 6)
    47 = @(7, 7)
    48 = @(7, 8)
    49 = @(9, 9)
    50 = @(6, 6)
    51 = @(7, 7)
    52 = @(4, 5)
    53 = @(7, 7)
    54 = @(7, 8)
    55 = @(8, 8)
    56 = @(7, 7)
    57 = @(3, 3)
    58 = @(9, 9)
    59 = @(9, 9)
    60 = @(6, 7)
    61 = @(9, 9)
    62 = @(9, 9)
    63 = @(9, 9)
    64 = @(9, 9)
    65 = @(9, 9)
    66 = @(8, 8)
    67 = @(8, 9)
    68 = @(9, 9)
    69 = @(9, 9)
    70 = @(6, 6)
    71 = @(6, 6)
    72 = @(6, 6)
    73 = @(9, 9)
    74 = @(3, 4)
    75 = @(9, 9)
    76 = @(8, 8)
    77 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item([int]$row, 9).Value = $vals[0]
    $ws.Cells.Item([int]$row, 10).Value = $vals[1]
}
